# Generate Report for handoff
#
# A new handoff was generated for the "fc355b15-..." source document's
# "993c6d4b-..." group row (row 4) in both the zh-cn and de-de worksheets.
# The handoff produced an identical target/handoff file (same content
# hash in the filename), so only the "Latest Handoff Datetime" (column D)
# timestamp advances for that row; every other cell keeps its text.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-17 09:55:03"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-17 09:55:14"
